$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.178.30"
$ws.Range("E2").Value = "'  -0.84%  "
$ws.Range("D3").Value = "'3.556.42"
$ws.Range("E3").Value = "'  -0.85%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "'  -0.13%  "
$ws.Range("D5").Value = "'617.51"
$ws.Range("E5").Value = "'  +5.09%  "
$ws.Range("D6").Value = "'185.94"
$ws.Range("E6").Value = "'  -0.09%  "
$ws.Range("E7").Value = "'  +0.82%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "'  -0.10%  "
$ws.Range("E9").Value = "'  +0.83%  "
$ws.Range("D10").Value = "'0.655"
$ws.Range("E10").Value = "'  +0.39%  "
$ws.Range("D11").Value = "'53.92"
$ws.Range("E11").Value = "'  -0.84%  "
$ws.Range("D12").Value = "'0.0000310"
$ws.Range("E12").Value = "'  -3.09%  "
$ws.Range("D13").Value = "'9.54"
$ws.Range("E13").Value = "'  -0.34%  "
$ws.Range("D14").Value = "'4.118.67"
$ws.Range("E14").Value = "'  -0.91%  "
$ws.Range("D15").Value = "'629.47"
$ws.Range("E15").Value = "'  +9.18%  "
$ws.Range("D16").Value = "'13.07"
$ws.Range("E16").Value = "'  +5.30%  "
$ws.Range("D17").Value = "'70.225.60"
$ws.Range("E17").Value = "'  -0.80%  "
$ws.Range("D18").Value = "'18.96"
$ws.Range("E18").Value = "'  -3.31%  "
$ws.Range("D19").Value = "'3.555.36"
$ws.Range("E19").Value = "'  -1.11%  "
$ws.Range("E20").Value = "'  -0.27%  "
$ws.Range("E21").Value = "'  -1.67%  "
$ws.Range("D22").Value = "'17.61"
$ws.Range("E22").Value = "'  -0.90%  "
$ws.Range("E23").Value = "'  +1.55%  "
$ws.Range("D24").Value = "'103.15"
$ws.Range("E24").Value = "'  +8.16%  "
$ws.Range("E25").Value = "'  +0.38%  "
$ws.Range("E26").Value = "'  +2.33%  "
$ws.Range("D27").Value = "'11.08"
$ws.Range("E27").Value = "'  -2.61%  "
$ws.Range("B28").Value = "'EthereumClassic"
$ws.Range("C28").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'34.48"
$ws.Range("E28").Value = "'  +6.77%  "
$ws.Range("B29").Value = "'Filecoin"
$ws.Range("C29").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "'9.56"
$ws.Range("E29").Value = "'  +4.35%  "
$ws.Range("D30").Value = "'7.08"
$ws.Range("E30").Value = "'  -3.31%  "
$ws.Range("D31").Value = "'12.31"
$ws.Range("E31").Value = "'  -0.84%  "
$ws.Range("E32").Value = "'  +0.22%  "
$ws.Range("E33").Value = "'  -1.75%  "
$ws.Range("E34").Value = "'  +20.35%  "
$ws.Range("D35").Value = "'3.27"
$ws.Range("E35").Value = "'  -2.82%  "
$ws.Range("D36").Value = "'531.79"
$ws.Range("E36").Value = "'  -5.55%  "
$ws.Range("E37").Value = "'  -3.92%  "
$ws.Range("E38").Value = "'  +0.00%  "
$ws.Range("D39").Value = "'37.30"
$ws.Range("E39").Value = "'  -1.14%  "
$ws.Range("D40").Value = "'3.578.84"
$ws.Range("E40").Value = "'  +5.56%  "
$ws.Range("D41").Value = "'0.0₃0783"
$ws.Range("E41").Value = "'  -1.41%  "
$ws.Range("D42").Value = "'3.54"
$ws.Range("E42").Value = "'  +4.45%  "
$ws.Range("D43").Value = "'0.138"
$ws.Range("E43").Value = "'  +1.62%  "
$ws.Range("E44").Value = "'  +2.59%  "
$ws.Range("D45").Value = "'2.96"
$ws.Range("E45").Value = "'  -0.65%  "
$ws.Range("D46").Value = "'0.143"
$ws.Range("E46").Value = "'  +4.03%  "
$ws.Range("D47").Value = "'3.38"
$ws.Range("E47").Value = "'  -5.21%  "
$ws.Range("D48").Value = "'9.19"
$ws.Range("E48").Value = "'  -1.55%  "
$ws.Range("E49").Value = "'  +0.26%  "
$ws.Range("E50").Value = "'  -1.59%  "
$ws.Range("D51").Value = "'134.34"
$ws.Range("E51").Value = "'  -2.17%  "
